$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Recipient data was refreshed: the two still-pending recipients moved up,
# two finished ones were dropped, and a brand-new recipient was appended.
# The trailing blank/placeholder rows also went away (8 rows -> 6 rows).
# ---------------------------------------------------------------------------

# Drop all existing hyperlinks first so they can be rebuilt cleanly against
# the new row layout (old rId1..rId3 would otherwise collide/duplicate).
$ws.Hyperlinks.Delete()

# Remove rows 4-8 outright (rather than just clearing contents) so the old
# "last row" thick-bottom-border flag and the now-unused blank rows 6-8
# disappear rather than lingering as empty styled cells.
$ws.Rows("4:8").Delete()

$notifUrl = "https://so-media-potral.vercel.app/passinterview/-O_75N9xzUIZDYnUcP5j"

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "`$NAME"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "`$NOTIFICATION_URL"

# --- Row 2: Đào Duy Thông (kept, notification url refreshed) -----------
$ws.Range("A2").Value = "Đào Duy Thông"
$ws.Range("B2").Value = "duythong.ptit@gmail.com"
$ws.Range("C2").Value = $notifUrl

# --- Row 3: Tống Ngọc Kiên (kept, notification url refreshed) ----------
$ws.Range("A3").Value = "Tống Ngọc Kiên"
$ws.Range("B3").Value = "duythong020703@gmail.com"
$ws.Range("C3").Value = $notifUrl

# --- Row 4: Vũ Thị Phương Thảo (now has her own hyperlinked row) -------
$ws.Range("A4").Value = "Vũ Thị Phương Thảo"
$ws.Range("B4").Value = "vtphth716@gmail.com"
$ws.Range("C4").Value = $notifUrl

# --- Row 5: Đào Dương Cẩm Tú (now has her own hyperlinked row) ---------
$ws.Range("A5").Value = "Đào Dương Cẩm Tú"
$ws.Range("B5").Value = "ddcamtus216@gmail.com"
$ws.Range("C5").Value = $notifUrl

# --- Row 6: Nguyễn Đoan Trang (brand-new recipient) ---------------------
$ws.Range("A6").Value = "Nguyễn Đoan Trang"
$ws.Range("B6").Value = "kimcotton124@gmail.com"
$ws.Range("C6").Value = $notifUrl

# Restore the normal (non-bold, non-boxed) row height/look for the
# recreated rows so they match rows 1-3.
$ws.Rows("4:6").RowHeight = 15.75

# Column A keeps the plain black Calibri look used throughout the sheet.
$ws.Range("A4:A6").Font.Name = "Calibri"
$ws.Range("A4:A6").Font.Size = 11
$ws.Range("A4:A6").Font.Color = 0

# --- Rebuild the mailto hyperlinks for the whole Email column ----------
# (Add in the same order the author filled them in: 2, 3, 4, 6, then 5.)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:duythong.ptit@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:duythong020703@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:vtphth716@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:kimcotton124@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:ddcamtus216@gmail.com")

# Hyperlinks.Add() mutates the cell's format; reapply the clean "Hyperlink"
# named style afterwards so every emailed cell looks like B3 did before.
$ws.Range("B3:B6").Style = "Hyperlink"

# B2 keeps its own distinctive blue-underline look (it was never the
# built-in "Hyperlink" style to begin with).
$ws.Range("B2").Font.Name = "Calibri"
$ws.Range("B2").Font.Size = 11
$ws.Range("B2").Font.Underline = 2
$ws.Range("B2").Font.Color = 12673797

# --- View bookkeeping: final selection left on D7 by the author --------
$ws.Range("D7").Select()
